$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2199.625  # was 1930.9131
$ws.Range("I116").Value = 1718.8  # was 1566.9375
$ws.Range("J116").Value = 3001  # was 2762.8572
$ws.Range("K116").Value = 1718.8  # was 1566.9375
$ws.Range("L116").Value = 3001  # was 2762.8572
$ws.Range("M116").Value = 1723.2  # was 1875.0625
$ws.Range("N116").Value = -9885  # was -9646.8572

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 25000540  # was 25641556
$ws.Range("I135").Value = 216.11429  # was 221.26471
$ws.Range("J135").Value = 200002800  # was 200002620
$ws.Range("K135").Value = 1945.02861  # was 1991.38239
$ws.Range("L135").Value = 1800025200  # was 1800023580
$ws.Range("M135").Value = 589.9713899999999  # was 543.61761
$ws.Range("N135").Value = -1800030270  # was -1800028650

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1184.5625  # was 1123.746
$ws.Range("I137").Value = 899.6842  # was 817.1556
$ws.Range("J137").Value = 1600.9231  # was 1890.2222
$ws.Range("K137").Value = 2699.0526  # was 2451.4668
$ws.Range("L137").Value = 4802.7693  # was 5670.6666
$ws.Range("M137").Value = -149.0526  # was 98.53319999999985
$ws.Range("N137").Value = -9902.7693  # was -10770.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 10000000  # was 1000
$ws.Range("I11").Value = 10000000  # was 1000
$ws.Range("K11").Value = 10000000  # was 1000
$ws.Range("M11").Value = -9999856  # was -856

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4451.0723  # was 4701.1846
$ws.Range("I32").Value = 4130.4746  # was 4402.7456
$ws.Range("K32").Value = 4130.4746  # was 4402.7456
$ws.Range("M32").Value = -3843.4746  # was -4115.7456

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 35715740  # was 37038508
$ws.Range("I61").Value = 43479544  # was 45455850
$ws.Range("J61").Value = 2242.8  # was 2222.8
$ws.Range("K61").Value = 43479544  # was 45455850
$ws.Range("L61").Value = 2242.8  # was 2222.8
$ws.Range("M61").Value = -43479332  # was -45455638
$ws.Range("N61").Value = -2666.8  # was -2646.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 966.24445  # was 981.38635
$ws.Range("I74").Value = 727.5526  # was 739.1081
$ws.Range("K74").Value = 727.5526  # was 739.1081
$ws.Range("M74").Value = 146.4474  # was 134.8919

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 966.24445  # was 981.38635
$ws.Range("I77").Value = 727.5526  # was 739.1081
$ws.Range("K77").Value = 3637.763  # was 3695.5405
$ws.Range("M77").Value = 730.2370000000001  # was 672.4594999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2493  # was 1868.4
$ws.Range("I110").Value = 2100  # was 1223.625
$ws.Range("K110").Value = 2100  # was 1223.625
$ws.Range("M110").Value = -55  # was 821.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1960.6  # was 2047.0714
$ws.Range("I122").Value = 2082.3076  # was 2193.3333
$ws.Range("K122").Value = 6246.9228  # was 6579.999899999999
$ws.Range("M122").Value = -3796.9228  # was -4129.999899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2334.96  # was 2069.8235
$ws.Range("J132").Value = 2200  # was 1490.909
$ws.Range("L132").Value = 6600  # was 4472.727000000001
$ws.Range("N132").Value = -11660  # was -9532.727000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 35715740  # was 37038508
$ws.Range("I136").Value = 43479544  # was 45455850
$ws.Range("J136").Value = 2242.8  # was 2222.8
$ws.Range("K136").Value = 130438632  # was 136367550
$ws.Range("L136").Value = 6728.400000000001  # was 6668.400000000001
$ws.Range("M136").Value = -130436082  # was -136365000
$ws.Range("N136").Value = -11828.4  # was -11768.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 440.5  # was 452.4
$ws.Range("J5").Value = 0  # was 500
$ws.Range("L5").Value = 0  # was 500
$ws.Range("N5").ClearContents()  # was -726

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4342.028  # was 3994.658
$ws.Range("I134").Value = 1100.3939  # was 1087.3636
$ws.Range("J134").Value = 40000  # was 23182.8
$ws.Range("K134").Value = 3301.1817  # was 3262.0908
$ws.Range("L134").Value = 120000  # was 69548.39999999999
$ws.Range("M134").Value = -766.1817000000001  # was -727.0907999999999
$ws.Range("N134").Value = -125070  # was -74618.39999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 76924344  # was 55556690
$ws.Range("I16").Value = 83334584  # was 58824640
$ws.Range("K16").Value = 83334584  # was 58824640
$ws.Range("M16").Value = -83334297  # was -58824353

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1899.4642  # was 1811.1666
$ws.Range("I31").Value = 1942.826  # was 1833.4
$ws.Range("K31").Value = 1942.826  # was 1833.4
$ws.Range("M31").Value = -1647.826  # was -1538.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1899.4642  # was 1811.1666
$ws.Range("I34").Value = 1942.826  # was 1833.4
$ws.Range("K34").Value = 1942.826  # was 1833.4
$ws.Range("M34").Value = -1740.826  # was -1631.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 76924344  # was 55556690
$ws.Range("I113").Value = 83334584  # was 58824640
$ws.Range("K113").Value = 83334584  # was 58824640
$ws.Range("M113").Value = -83332414  # was -58822470

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 22728256  # was 20834418
$ws.Range("I134").Value = 1005.8  # was 1032.421
$ws.Range("J134").Value = 250000750  # was 100001280
$ws.Range("K134").Value = 3017.4  # was 3097.263
$ws.Range("L134").Value = 750002250  # was 300003840
$ws.Range("M134").Value = -482.3999999999996  # was -562.2629999999999
$ws.Range("N134").Value = -750007320  # was -300008910

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 5700  # was 4760
$ws.Range("I5").Value = 5700  # was 4760
$ws.Range("K5").Value = 17100  # was 14280
$ws.Range("M5").Value = -16988  # was -14168

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 1959.2727  # was 1442
$ws.Range("I32").Value = 1067.3334  # was 1025.5
$ws.Range("J32").Value = 2293.75  # was 2275
$ws.Range("K32").Value = 3202.0002  # was 3076.5
$ws.Range("L32").Value = 6881.25  # was 6825
$ws.Range("M32").Value = -2919.0002  # was -2793.5
$ws.Range("N32").Value = -7447.25  # was -7391

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 3560  # was 2842.8572
$ws.Range("I80").Value = 1033.3334  # was 1020
$ws.Range("J80").Value = 4642.857  # was 3855.5557
$ws.Range("K80").Value = 3100.0002  # was 3060
$ws.Range("L80").Value = 13928.571  # was 11566.6671
$ws.Range("M80").Value = -2164.0002  # was -2124
$ws.Range("N80").Value = -15800.571  # was -13438.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 3560  # was 2842.8572
$ws.Range("I83").Value = 1033.3334  # was 1020
$ws.Range("J83").Value = 4642.857  # was 3855.5557
$ws.Range("K83").Value = 9300.000599999999  # was 9180
$ws.Range("L83").Value = 41785.713  # was 34700.0013
$ws.Range("M83").Value = -4620.000599999999  # was -4500
$ws.Range("N83").Value = -51145.713  # was -44060.0013

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 265.52  # was 240.59259
$ws.Range("I92").Value = 272.54544  # was 221.14285
$ws.Range("J92").Value = 260  # was 261.53845
$ws.Range("K92").Value = 817.63632  # was 663.4285500000001
$ws.Range("L92").Value = 780  # was 784.61535
$ws.Range("M92").Value = 430.36368  # was 584.5714499999999
$ws.Range("N92").Value = -3276  # was -3280.61535

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 5639.4736  # was 5365
$ws.Range("J107").Value = 7008.7334  # was 6580.0625
$ws.Range("L107").Value = 21026.2002  # was 19740.1875
$ws.Range("N107").Value = -24866.2002  # was -23580.1875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 730.625  # was 730.95654
$ws.Range("J113").Value = 730.625  # was 730.95654
$ws.Range("L113").Value = 2191.875  # was 2192.86962
$ws.Range("N113").Value = -6531.875  # was -6532.869619999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 18183058  # was 15152653
$ws.Range("J131").Value = 1364.234  # was 1224.5518
$ws.Range("L131").Value = 4092.702  # was 3673.6554
$ws.Range("N131").Value = -14172.702  # was -13753.6554

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 3322.3333  # was 3535.3447
$ws.Range("I134").Value = 1378.9286  # was 1539.5834
$ws.Range("J134").Value = 5022.8125  # was 4944.1177
$ws.Range("K134").Value = 4136.7858  # was 4618.7502
$ws.Range("L134").Value = 15068.4375  # was 14832.3531
$ws.Range("M134").Value = 933.2142000000003  # was 451.2497999999996
$ws.Range("N134").Value = -25208.4375  # was -24972.3531

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 5700  # was 4760
$ws.Range("I135").Value = 5700  # was 4760
$ws.Range("K135").Value = 51300  # was 42840
$ws.Range("M135").Value = -48765  # was -40305

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1415.25  # was 1431.2667
$ws.Range("I113").Value = 1349.4546  # was 1365.6
$ws.Range("J113").Value = 1560  # was 1562.6
$ws.Range("K113").Value = 1349.4546  # was 1365.6
$ws.Range("L113").Value = 1560  # was 1562.6
$ws.Range("M113").Value = 820.5454  # was 804.4000000000001
$ws.Range("N113").Value = -5900  # was -5902.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2313.5789  # was 2178.9524
$ws.Range("I122").Value = 1730.6  # was 1632.8823
$ws.Range("K122").Value = 5191.799999999999  # was 4898.6469
$ws.Range("M122").Value = -2741.799999999999  # was -2448.6469

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2038.1111  # was 2218
$ws.Range("I132").Value = 1442.5714  # was 1626
$ws.Range("J132").Value = 4122.5  # was 3500.6667
$ws.Range("K132").Value = 4327.7142  # was 4878
$ws.Range("L132").Value = 12367.5  # was 10502.0001
$ws.Range("M132").Value = -1797.7142  # was -2348
$ws.Range("N132").Value = -17427.5  # was -15562.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 500  # was 0
$ws.Range("J14").Value = 500  # was 0
$ws.Range("L14").Value = 500  # was 0
$ws.Range("N14").Value = -844  # was None

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 22887.723  # was 22431.688
$ws.Range("I132").Value = 1237.4667  # was 1240.7333
$ws.Range("J132").Value = 61094.06  # was 57749.945
$ws.Range("K132").Value = 3712.4001  # was 3722.199900000001
$ws.Range("L132").Value = 183282.18  # was 173249.835
$ws.Range("M132").Value = -1182.4001  # was -1192.199900000001
$ws.Range("N132").Value = -188342.18  # was -178309.835

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3149.3447  # was 2832.7273
$ws.Range("I132").Value = 3881.9524  # was 3598.739
$ws.Range("J132").Value = 1226.25  # was 1070.9
$ws.Range("K132").Value = 11645.8572  # was 10796.217
$ws.Range("L132").Value = 3678.75  # was 3212.7
$ws.Range("M132").Value = -9115.8572  # was -8266.217000000001
$ws.Range("N132").Value = -8738.75  # was -8272.700000000001
